$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Move "PowerPlants" sheet tab so it comes right before "FuelsExisting"
# ------------------------------------------------------------------
$ppWs = $wb.Worksheets.Item("PowerPlants")
$feWs = $wb.Worksheets.Item("FuelsExisting")
$ppWs.Move($feWs)

# ------------------------------------------------------------------
# 2. Add the new "CapacityCredit" column (column I) to PowerPlants
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PowerPlants")

$ws.Cells.Item(1, 9).Value = "CapacityCredit"
$ws.Cells.Item(2, 9).Value = "[fr]"

$ws.Cells.Item(3, 9).Value = 0.9
$ws.Cells.Item(4, 9).Value = 0.91
$ws.Cells.Item(5, 9).Value = 0.95
$ws.Cells.Item(6, 9).Value = 0.91
$ws.Cells.Item(7, 9).Value = 0.91
$ws.Cells.Item(8, 9).Value = 0.91
$ws.Cells.Item(9, 9).Value = 0.94
$ws.Cells.Item(10, 9).Value = 0.95
$ws.Cells.Item(11, 9).Value = 0.95
$ws.Cells.Item(12, 9).Value = 0.29
$ws.Cells.Item(13, 9).Value = 0.36
$ws.Cells.Item(14, 9).Value = 0.75
$ws.Cells.Item(15, 9).Value = 0.95
$ws.Cells.Item(16, 9).Value = 0.9
$ws.Cells.Item(17, 9).Value = 0.95
$ws.Cells.Item(18, 9).Value = 0.95
$ws.Cells.Item(19, 9).Value = 0.95
$ws.Cells.Item(20, 9).Value = 0.91
$ws.Cells.Item(21, 9).Value = 0.29
$ws.Cells.Item(22, 9).Value = 0.36
$ws.Cells.Item(23, 9).Value = 0.75
$ws.Cells.Item(24, 9).Value = 0.95
$ws.Cells.Item(25, 9).Value = 0.95
$ws.Cells.Item(26, 9).Value = 0.95
$ws.Cells.Item(27, 9).Value = 0.29
$ws.Cells.Item(28, 9).Value = 0.36

# ------------------------------------------------------------------
# 3. Split the existing "Y" conditional formatting so that column E has
#    its own rule(s), separate from D and F:H (same dxf as before for
#    E, fresh duplicated dxf for D/F:H).
# ------------------------------------------------------------------
$ws.Range("D3:H15").FormatConditions.Delete()
$ws.Range("D16:H18").FormatConditions.Delete()

$rule1 = $ws.Range("D3:D15,F3:H15").FormatConditions.Add(1, 3, '"Y"')
$rule1.Interior.Color = $ws.Range("D3").Interior.Color
$rule2 = $ws.Range("D19:D28,F19:H28").FormatConditions.Add(1, 3, '"Y"')
$rule2.Interior.Color = $ws.Range("D3").Interior.Color

$rule3 = $ws.Range("D16:D18,F16:H18").FormatConditions.Add(1, 3, '"Y"')
$rule3.Interior.Color = $ws.Range("D3").Interior.Color

$rule4 = $ws.Range("E3:E15,E19:E28").FormatConditions.Add(1, 3, '"Y"')
$rule4.Interior.Color = $ws.Range("D3").Interior.Color

$rule5 = $ws.Range("E16:E18").FormatConditions.Add(1, 3, '"Y"')
$rule5.Interior.Color = $ws.Range("D3").Interior.Color

# ------------------------------------------------------------------
# 4. Reset column E's cell style on the data rows (the border/fill
#    "applied" flags get cleared when Excel re-lays the formatting out)
# ------------------------------------------------------------------
$ws.Range("E3:E28").ClearFormats()

Write-Output "done"
